# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The "K" column (G) holds recalculated values for rows 2-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 1
    9  = 0
    10 = 3
    11 = 0
    12 = 1
    13 = 1
    14 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
